# Hortaliza, Vega Modelo de Temuco - Espinaca
# Weekly update: a new week's record is inserted at the top of the
# "Espinaca" data block (row 74), pushing the existing rows 74-120 down
# to rows 75-121 and extending the used range from A1:R120 to A1:R121.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 74; this shifts every
# subsequent row (old 74..120) down by one (new 75..121), which is
# exactly what the target diff shows (each old row's content reappears
# one row lower, and a brand-new row 74 appears with fresh data).
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with this week's record. Columns
# A, B, C, E, F, G, H, I, N, O, Q, R repeat the same constant values used
# throughout this block (market/region/category metadata), while D
# (fecha), J (volumen), K/L/M (precios), and P (precio $/Kg) carry the
# new week's figures.
$ws.Cells.Item(74, 1).Value = 10
$ws.Cells.Item(74, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(74, 3).Value = "La Araucanía"
$ws.Cells.Item(74, 4).Value = 44582
$ws.Cells.Item(74, 5).Value = 9
$ws.Cells.Item(74, 6).Value = 100112012
$ws.Cells.Item(74, 7).Value = "Espinaca"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 30
$ws.Cells.Item(74, 11).Value = 14000
$ws.Cells.Item(74, 12).Value = 14000
$ws.Cells.Item(74, 13).Value = 14000
$ws.Cells.Item(74, 14).Value = "`$/docena de atados"
$ws.Cells.Item(74, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(74, 16).Value = 4667
$ws.Cells.Item(74, 17).Value = 3
$ws.Cells.Item(74, 18).Value = "Hortaliza"
